$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Roll No (A2): was text "22101B0022" -> becomes numeric 1234
$ws.Range("A2").Value = 1234

# Email (D2): "mohitkoli@gmail.com" -> "mohitkoli1234@gmail.com"
$ws.Range("D2").Value = "mohitkoli1234@gmail.com"

# Password (E2): was text "mohit" -> becomes numeric 123
$ws.Range("E2").Value = 123

# Move the active cell selection, matching the saved cursor position
$ws.Range("D10").Select() | Out-Null
